$d = $word.ActiveDocument

# The court appearance / license-suspension date moved from July 04, 2022
# to July 06, 2022. This phrase occurs three times in the document
# (the narrative paragraph, the chart cell, and the license-suspension
# sentence), so replace every occurrence in one pass.
$d.Content.Find.Execute("July 04, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "July 06, 2022", 2) | Out-Null

# The community-control reporting deadline moved from September 02, 2022
# to September 04, 2022.
$d.Content.Find.Execute("September 02, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 04, 2022", 2) | Out-Null
